# Apply the edit described by the diff:
# - 3 new data rows ("Región de O'Higgins" origin, fecha 2023-01-20) are inserted
#   right before the existing row 540 (which holds "Sandia" data for Femacal de
#   La Calera), shifting all subsequent rows (540-624) down by 3 (to 543-627).
# - The worksheet dimension grows from A1:R624 to A1:R627.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 540 (this shifts existing rows 540:624 down to 543:627)
$insertRange = $ws.Range("A540:R542")
$insertRange.Insert()

# Fill in the 3 newly inserted rows with the new data from the diff.
# Row 540
$ws.Cells.Item(540, 1).Value = 3
$ws.Cells.Item(540, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(540, 3).Value = "Coquimbo"
$ws.Cells.Item(540, 4).Value = [datetime]"2023-01-20"
$ws.Cells.Item(540, 5).Value = 5
$ws.Cells.Item(540, 6).Value = 100112028
$ws.Cells.Item(540, 7).Value = "Sandia"
$ws.Cells.Item(540, 8).Value = "Sin especificar"
$ws.Cells.Item(540, 9).Value = "Extra"
$ws.Cells.Item(540, 10).Value = 740
$ws.Cells.Item(540, 11).Value = 3000
$ws.Cells.Item(540, 12).Value = 3500
$ws.Cells.Item(540, 13).Value = 3324
$ws.Cells.Item(540, 14).Value = "$/unidad"
$ws.Cells.Item(540, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(540, 16).Value = 3324
$ws.Cells.Item(540, 17).Value = 1
$ws.Cells.Item(540, 18).Value = "Hortaliza"

# Row 541
$ws.Cells.Item(541, 1).Value = 3
$ws.Cells.Item(541, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(541, 3).Value = "Coquimbo"
$ws.Cells.Item(541, 4).Value = [datetime]"2023-01-20"
$ws.Cells.Item(541, 5).Value = 5
$ws.Cells.Item(541, 6).Value = 100112028
$ws.Cells.Item(541, 7).Value = "Sandia"
$ws.Cells.Item(541, 8).Value = "Sin especificar"
$ws.Cells.Item(541, 9).Value = "Primera"
$ws.Cells.Item(541, 10).Value = 830
$ws.Cells.Item(541, 11).Value = 2300
$ws.Cells.Item(541, 12).Value = 2500
$ws.Cells.Item(541, 13).Value = 2408
$ws.Cells.Item(541, 14).Value = "$/unidad"
$ws.Cells.Item(541, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(541, 16).Value = 2408
$ws.Cells.Item(541, 17).Value = 1
$ws.Cells.Item(541, 18).Value = "Hortaliza"

# Row 542
$ws.Cells.Item(542, 1).Value = 3
$ws.Cells.Item(542, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(542, 3).Value = "Coquimbo"
$ws.Cells.Item(542, 4).Value = [datetime]"2023-01-20"
$ws.Cells.Item(542, 5).Value = 5
$ws.Cells.Item(542, 6).Value = 100112028
$ws.Cells.Item(542, 7).Value = "Sandia"
$ws.Cells.Item(542, 8).Value = "Sin especificar"
$ws.Cells.Item(542, 9).Value = "Segunda"
$ws.Cells.Item(542, 10).Value = 500
$ws.Cells.Item(542, 11).Value = 1800
$ws.Cells.Item(542, 12).Value = 1800
$ws.Cells.Item(542, 13).Value = 1800
$ws.Cells.Item(542, 14).Value = "$/unidad"
$ws.Cells.Item(542, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(542, 16).Value = 1800
$ws.Cells.Item(542, 17).Value = 1
$ws.Cells.Item(542, 18).Value = "Hortaliza"
